$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "n_successful" column (Q) entirely; remaining columns shift dimension to A1:P6
$ws.Range("Q1:Q6").EntireColumn.Delete()

# Update header row (M1:P1) to reflect the new column structure
$ws.Range("M1").Value = "faithfulness"
$ws.Range("N1").Value = "overall"
$ws.Range("O1").Value = "n_cases"
$ws.Range("P1").Value = "n_successful"

# Update data rows 2-6 with refreshed evaluation numbers
$ws.Range("C2").Value = 0.418
$ws.Range("D2").Value = 0.465
$ws.Range("E2").Value = 0.13
$ws.Range("F2").Value = 0.276
$ws.Range("G2").Value = 0.095
$ws.Range("H2").Value = 0.129
$ws.Range("I2").Value = 0.024
$ws.Range("J2").Value = 0.125
$ws.Range("K2").Value = -0.057
$ws.Range("L2").Value = 0.51
$ws.Range("M2").Value = 0.51
$ws.Range("N2").Value = 0.255
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 24

$ws.Range("C3").Value = 0.408
$ws.Range("D3").Value = 0.472
$ws.Range("E3").Value = 0.144
$ws.Range("F3").Value = 0.306
$ws.Range("G3").Value = 0.093
$ws.Range("H3").Value = 0.134
$ws.Range("I3").Value = 0.04
$ws.Range("J3").Value = 0.147
$ws.Range("K3").Value = -0.075
$ws.Range("L3").Value = 0.528
$ws.Range("M3").Value = 0.528
$ws.Range("N3").Value = 0.256
$ws.Range("O3").Value = 25
$ws.Range("P3").Value = 25

$ws.Range("C4").Value = 0.358
$ws.Range("D4").Value = 0.408
$ws.Range("E4").Value = 0.119
$ws.Range("F4").Value = 0.255
$ws.Range("G4").Value = 0.095
$ws.Range("H4").Value = 0.13
$ws.Range("I4").Value = 0.012
$ws.Range("J4").Value = 0.104
$ws.Range("K4").Value = -0.003
$ws.Range("L4").Value = 0.58
$ws.Range("M4").Value = 0.58
$ws.Range("N4").Value = 0.259
$ws.Range("O4").Value = 25
$ws.Range("P4").Value = 25

$ws.Range("C5").Value = 0.332
$ws.Range("D5").Value = 0.358
$ws.Range("E5").Value = 0.145
$ws.Range("F5").Value = 0.304
$ws.Range("G5").Value = 0.11
$ws.Range("H5").Value = 0.141
$ws.Range("I5").Value = 0.03
$ws.Range("J5").Value = 0.14
$ws.Range("K5").Value = -0.021
$ws.Range("L5").Value = 0.501
$ws.Range("M5").Value = 0.501
$ws.Range("N5").Value = 0.239
$ws.Range("O5").Value = 25
$ws.Range("P5").Value = 23

$ws.Range("C6").Value = 0.317
$ws.Range("D6").Value = 0.357
$ws.Range("E6").Value = 0.192
$ws.Range("F6").Value = 0.394
$ws.Range("G6").Value = 0.128
$ws.Range("H6").Value = 0.151
$ws.Range("I6").Value = 0.063
$ws.Range("J6").Value = 0.222
$ws.Range("K6").Value = -0.003
$ws.Range("L6").Value = 0.484
$ws.Range("M6").Value = 0.484
$ws.Range("N6").Value = 0.249
$ws.Range("O6").Value = 25
$ws.Range("P6").Value = 22

